$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.855.69"
$ws.Range("E2").Value = "  +1.49%  "

$ws.Range("D3").Value = "3.547.88"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'614.11"
$ws.Range("E5").Value = "  +1.57%  "

$ws.Range("D6").Value = "'152.74"
$ws.Range("E6").Value = "  -0.94%  "

$ws.Range("D7").Value = "3.547.25"
$ws.Range("E7").Value = "  +0.63%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("E11").Value = "  +3.41%  "

$ws.Range("D12").Value = "'0.427"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("D14").Value = "4.143.57"
$ws.Range("E14").Value = "  +0.48%  "

$ws.Range("D15").Value = "'32.01"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").Value = "3.544.66"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "67.630.23"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").Value = "'0.116"
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("D20").Value = "'15.24"
$ws.Range("E20").Value = "  -1.33%  "

$ws.Range("D21").Value = "'9.71"
$ws.Range("E21").Value = "  +3.73%  "

$ws.Range("D22").Value = "'447.33"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("E23").Value = "  -2.11%  "

$ws.Range("D24").Value = "'77.03"
$ws.Range("E24").Value = "  -2.38%  "

$ws.Range("D25").Value = "'0.0000130"
$ws.Range("E25").Value = "  +5.83%  "

$ws.Range("D26").Value = "3.688.42"
$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").Value = "'10.22"
$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("E29").Value = "  +4.70%  "

$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("D31").Value = "'1.62"
$ws.Range("E31").Value = "  -3.29%  "

$ws.Range("E32").Value = "  +6.37%  "

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").Value = "'25.82"
$ws.Range("E34").Value = "  -0.32%  "

$ws.Range("D35").Value = "'6.22"
$ws.Range("E35").Value = "  +0.51%  "

$ws.Range("D36").Value = "3.533.44"
$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("D38").Value = "'8.05"
$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.13%  "

$ws.Range("D41").Value = "'2.22"
$ws.Range("E41").Value = "  +3.60%  "

$ws.Range("D42").Value = "'176.21"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("D43").Value = "'0.0895"
$ws.Range("E43").Value = "  +2.25%  "

$ws.Range("E44").Value = "  -3.34%  "

$ws.Range("E45").Value = "  -0.49%  "

$ws.Range("D46").Value = "'28.98"
$ws.Range("E46").Value = "  +2.64%  "

$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("E49").Value = "  +3.99%  "

$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "'0.253"
$ws.Range("E51").Value = "  +0.96%  "
